$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new F column ("dSF") value
$changes = @{
    5  = 8
    6  = -3
    14 = -8
    17 = -1
    25 = 3
    32 = -8
    33 = 0
    41 = 2
    45 = 6
    47 = -2
    49 = -2
    51 = -3
    55 = 8
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
